# events_template.xlsx refactor:
#  - remove the "location" column (free-text place name)
#  - add "latitude"/"longitude" numeric columns instead
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "location" was column C; dropping it shifts startDateTime/endDateTime left
# and the former selection/columns recompute automatically.
$ws.Columns("C").Delete()

# New trailing columns for geo-coordinates replacing the free text location.
$ws.Range("E1").Value = "latitude"
$ws.Range("F1").Value = "longitude"

$ws.Range("E2").Value = 4.695014
$ws.Range("F2").Value = -74.116591

$ws.Range("E3").Value = 4.695014
$ws.Range("F3").Value = -74.116591

$ws.Range("E4").Value = 4.695014
$ws.Range("F4").Value = -74.116591

# Cosmetic re-layout that follows from the column refactor.
$ws.Columns("A:B").ColumnWidth = 21.75
$ws.Columns("C").ColumnWidth = 21.584
$ws.Columns("D").ColumnWidth = 22.417
$ws.Columns("E").ColumnWidth = 23.917

$ws.Rows("1").RowHeight = 13.8

[void]$ws.Range("H3").Select()
